# "Updated Batch addon Step"
# - Msg sheet becomes the active/selected sheet (tab), replacing Login.
# - Msg!A11 gets a new validation message ("Desc Strt Spl"), added as a
#   brand-new shared string.
# - Selection on the Msg sheet moves from C11 to B11.

$wb = $excel.ActiveWorkbook

$wsMsg = $wb.Worksheets.Item("Msg")

# Make "Msg" the active sheet (this flips tabSelected off "Login" and onto
# "Msg", and updates the workbook's activeTab index).
$wsMsg.Activate()

# New test-scenario message text for the "Desc Strt Spl" case.
$wsMsg.Range("A11").Value = "Desc Strt Spl"

# Move the in-sheet selection to B11.
[void]$wsMsg.Range("B11").Select()
